$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates for rows 3-15 (row 5 unchanged),
# applied identically to both "展览" and "全部类型" sheets.
$updates = @{
    3  = 37
    4  = 139
    6  = 474
    7  = 1326
    8  = 426
    9  = 90
    10 = 158
    11 = 110
    12 = 161
    13 = 94
    14 = 140
    15 = 130
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
